$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New "disponible" stock values for column F (rows 2-40), replacing the
# placeholder text "1" with real quantities, formatted as numbers.
$values = @(13,3,9,23,36,9,6,13,34,27,25,25,41,26,11,17,29,26,21,8,18,9,32,2,8,46,10,18,7,41,0,30,13,5,40,62,15,10,40)

$ws.Range("F2:F40").NumberFormat = "#,##0"

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}

$ws.Range("E2").Select()
